$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The survey data gets an extra quarter of data: duplicate the last data
# row (row 13) into a brand-new row 14 (same values & formatting), then
# correct the "Ended date" on the original row (row 13) and give the new
# row its own (different) "Ended date".

# 1) Duplicate row 13 (values + formatting) into row 14.
$ws.Rows.Item(13).Copy()
$ws.Rows.Item(14).PasteSpecial(-4104)

# 2) The "Ended date" cell (column I) on the new row should look like the
#    other "General"-formatted date-as-text cells (e.g. I12), not like the
#    "m/d/yyyy"-formatted one it inherited from row 13.
$ws.Range("I12").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = "'2022-07-10"

# 3) Correct row 13's "Ended date" value (kept in its original date-looking
#    text format).
$ws.Range("I13").Value = "'2022-05-10"

# 4) Reflect the author's final cursor position in the saved file.
$ws.Range("I12").Select() | Out-Null
